# Scheduled runner update: refresh currentAveragePrice / profit figures
# across the ALC, ARM, BSM, CRP, CUL, GSM, LTW and WVR leve-profit sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H94").Value = 9999.667
$ws.Range("I94").Value = 9999.5
$ws.Range("K94").Value = 9999.5
$ws.Range("M94").Value = -9548.5

$ws.Range("H100").Value = 2578.8276
$ws.Range("I100").Value = 2140.682
$ws.Range("J100").Value = 3955.8572
$ws.Range("K100").Value = 2140.682
$ws.Range("L100").Value = 3955.8572
$ws.Range("M100").Value = -1599.682
$ws.Range("N100").Value = -5037.8572

$ws.Range("H127").Value = 2242
$ws.Range("I127").Value = 2249
$ws.Range("J127").Value = 2200
$ws.Range("K127").Value = 6747
$ws.Range("L127").Value = 6600
$ws.Range("M127").Value = -1787
$ws.Range("N127").Value = -16520

$ws.Range("H129").Value = 2378.4814
$ws.Range("I129").Value = 3397
$ws.Range("J129").Value = 1949.6316
$ws.Range("K129").Value = 10191
$ws.Range("L129").Value = 5848.8948
$ws.Range("M129").Value = -5191
$ws.Range("N129").Value = -15848.8948

$ws.Range("H131").Value = 4963.8
$ws.Range("I131").Value = 781
$ws.Range("J131").Value = 6756.4287
$ws.Range("K131").Value = 2343
$ws.Range("L131").Value = 20269.2861
$ws.Range("M131").Value = 2697
$ws.Range("N131").Value = -30349.2861

$ws.Range("H141").Value = 6730.6216
$ws.Range("I141").Value = 6121.6177
$ws.Range("K141").Value = 18364.8531
$ws.Range("M141").Value = -13184.8531

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 9796.395
$ws.Range("I32").Value = 7039.1455
$ws.Range("K32").Value = 7039.1455
$ws.Range("M32").Value = -6752.1455

$ws.Range("H45").Value = 316268.28
$ws.Range("I45").Value = 835592.75
$ws.Range("J45").Value = 4673.6
$ws.Range("K45").Value = 835592.75
$ws.Range("L45").Value = 4673.6
$ws.Range("M45").Value = -835215.75
$ws.Range("N45").Value = -5427.6

$ws.Range("H61").Value = 3877.2207
$ws.Range("I61").Value = 3749.1404
$ws.Range("K61").Value = 3749.1404
$ws.Range("M61").Value = -3537.1404

$ws.Range("H132").Value = 27623.904
$ws.Range("I132").Value = 30562.432
$ws.Range("K132").Value = 91687.296
$ws.Range("M132").Value = -89157.296

$ws.Range("H136").Value = 3877.2207
$ws.Range("I136").Value = 3749.1404
$ws.Range("K136").Value = 11247.4212
$ws.Range("M136").Value = -8697.4212

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H125").Value = 92189.5
$ws.Range("J125").Value = 92189.5
$ws.Range("L125").Value = 92189.5
$ws.Range("N125").Value = -102029.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6451.2354
$ws.Range("I31").Value = 4316.5557
$ws.Range("J31").Value = 8852.75
$ws.Range("K31").Value = 4316.5557
$ws.Range("L31").Value = 8852.75
$ws.Range("M31").Value = -4021.5557
$ws.Range("N31").Value = -9442.75

$ws.Range("H34").Value = 6451.2354
$ws.Range("I34").Value = 4316.5557
$ws.Range("J34").Value = 8852.75
$ws.Range("K34").Value = 4316.5557
$ws.Range("L34").Value = 8852.75
$ws.Range("M34").Value = -4114.5557
$ws.Range("N34").Value = -9256.75

$ws.Range("H58").Value = 74297.5
$ws.Range("I58").Value = 86279.25
$ws.Range("K58").Value = 86279.25
$ws.Range("M58").Value = -86076.25

$ws.Range("H99").Value = 5999.2
$ws.Range("J99").Value = 5999.2
$ws.Range("L99").Value = 5999.2
$ws.Range("N99").Value = -8995.2

$ws.Range("H122").Value = 1800.3572
$ws.Range("I122").Value = 1800.3572
$ws.Range("K122").Value = 5401.071599999999
$ws.Range("M122").Value = -2951.071599999999

$ws.Range("H126").Value = 5999.2
$ws.Range("J126").Value = 5999.2
$ws.Range("L126").Value = 17997.6
$ws.Range("N126").Value = -22937.6

$ws.Range("H132").Value = 2684.2856
$ws.Range("I132").Value = 2734.8235
$ws.Range("K132").Value = 8204.4705
$ws.Range("M132").Value = -5674.470499999999

$ws.Range("H134").Value = 40186.777
$ws.Range("I134").Value = 50747.906
$ws.Range("K134").Value = 152243.718
$ws.Range("M134").Value = -149708.718

$ws.Range("H136").Value = 74297.5
$ws.Range("I136").Value = 86279.25
$ws.Range("K136").Value = 258837.75
$ws.Range("M136").Value = -256287.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 1565387.8
$ws.Range("I4").Value = 1657175.2
$ws.Range("K4").Value = 4971525.6
$ws.Range("M4").Value = -4971413.6

$ws.Range("H37").Value = 63598.285
$ws.Range("J37").Value = 63598.285
$ws.Range("L37").Value = 190794.855
$ws.Range("N37").Value = -191018.855

$ws.Range("H129").Value = 532768.94
$ws.Range("I129").Value = 1575.4445
$ws.Range("J129").Value = 1215732
$ws.Range("K129").Value = 4726.333500000001
$ws.Range("L129").Value = 3647196
$ws.Range("M129").Value = 273.6664999999994
$ws.Range("N129").Value = -3657196

$ws.Range("H131").Value = 8357333
$ws.Range("J131").Value = 11142222
$ws.Range("L131").Value = 33426666
$ws.Range("N131").Value = -33436746

$ws.Range("H140").Value = 233420.47
$ws.Range("I140").Value = 233420.47
$ws.Range("K140").Value = 700261.41
$ws.Range("M140").Value = -695081.41

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 3136.4285
$ws.Range("I102").Value = 3136.4285
$ws.Range("K102").Value = 3136.4285
$ws.Range("M102").Value = -1514.4285

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2853.4
$ws.Range("J22").Value = 3594.6667
$ws.Range("L22").Value = 3594.6667
$ws.Range("N22").Value = -4184.6667

$ws.Range("H27").Value = 2853.4
$ws.Range("J27").Value = 3594.6667
$ws.Range("L27").Value = 3594.6667
$ws.Range("N27").Value = -3808.6667

$ws.Range("H55").Value = 911.3571
$ws.Range("I55").Value = 419.7143
$ws.Range("J55").Value = 1403
$ws.Range("K55").Value = 419.7143
$ws.Range("L55").Value = 1403
$ws.Range("M55").Value = -246.7143
$ws.Range("N55").Value = -1749

$ws.Range("H82").Value = 3032.2222
$ws.Range("I82").Value = 1694.5555
$ws.Range("J82").Value = 4369.8887
$ws.Range("K82").Value = 1694.5555
$ws.Range("L82").Value = 4369.8887
$ws.Range("M82").Value = -1333.5555
$ws.Range("N82").Value = -5091.8887

$ws.Range("H85").Value = 3032.2222
$ws.Range("I85").Value = 1694.5555
$ws.Range("J85").Value = 4369.8887
$ws.Range("K85").Value = 1694.5555
$ws.Range("L85").Value = 4369.8887
$ws.Range("M85").Value = -446.5554999999999
$ws.Range("N85").Value = -6865.8887

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2884.111
$ws.Range("I122").Value = 2766.6155
$ws.Range("K122").Value = 8299.8465
$ws.Range("M122").Value = -5849.8465
